$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Two Sum"
$ws.Range("C2").Value = "Given an array of integers nums and an integer target, return indices of the two numbers such that they add up to target."
$ws.Range("D2").Value = "Hash Map"
$ws.Range("E2").Value = "Arrays"
$ws.Range("F2").Value = "Easy"
